# Update cryptos list values (prices and 1h volume percentages) to reflect latest snapshot.
# Also fixes row-ordering swaps for rows 43-46 (BinanceUSD/FraxShare and Cronos/ARBITRUM).
#
# Price-column values that look like plain numbers (e.g. "1.00", "6.70") must stay
# TEXT cells (matching the feed's original inline-string formatting), instead of being
# silently coerced to numeric by the usual text-looks-like-a-number auto-detect. We force
# that by stamping the cell as Text before the write, then restore the default "Normal"
# style afterwards so no stray number-format is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.946.74'
$ws.Range('E2').Value = '  +1.54%  '

$ws.Range('D3').Value = '2.374.15'
$ws.Range('E3').Value = '  +0.79%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.692'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.45%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.64'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.46%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '77.03'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +7.67%  '

$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.636'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +28.28%  '

$ws.Range('E10').Value = '  +5.37%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '57.55'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.30%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '33.06'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +20.72%  '

$ws.Range('E13').Value = '  +19.62%  '

$ws.Range('E14').Value = '  +2.19%  '

$ws.Range('D15').Value = '2.725.64'
$ws.Range('E15').Value = '  +0.55%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '17.04'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.41%  '

$ws.Range('E17').Value = '  +7.11%  '

$ws.Range('D18').Value = '2.363.34'
$ws.Range('E18').Value = '  +0.15%  '

$ws.Range('D19').Value = '45.385.61'
$ws.Range('E19').Value = '  +4.86%  '

$ws.Range('E20').Value = '  +2.18%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.70'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.44%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '78.12'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.72%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '258.62'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.29%  '

$ws.Range('E24').Value = '  +0.02%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.56'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.76%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.25'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +11.47%  '

$ws.Range('E27').Value = '  -2.14%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.78'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +17.17%  '

$ws.Range('E29').Value = '  +1.72%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '23.31'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.12%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '176.01'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.10%  '

$ws.Range('E32').Value = '  -0.29%  '

$ws.Range('E33').Value = '  +5.72%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.36'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.89%  '

$ws.Range('E35').Value = '  +9.03%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.39'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.70%  '

$ws.Range('E37').Value = '  +2.63%  '

$ws.Range('E38').Value = '  +1.88%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.52'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.49%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0276'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.40%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '19.10'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.64%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.203'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +19.08%  '

$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.95'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.43%  '

$ws.Range('B44').Value = 'BinanceUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.10%  '

$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.21'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.46%  '

$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.101'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.55%  '

$ws.Range('E47').Value = '  +5.13%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.53'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +14.26%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '103.05'
$ws.Range('D49').Style = 'Normal'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.51'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.11%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '54.78'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.57%  '

